$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.276.48"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.830.67"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "'235.60"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'0.6027"
$ws.Range("E6").Value = "  -3.81%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -4.97%  "
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("E10").Value = "  -5.52%  "
$ws.Range("D11").Value = "'0.07666"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.838.79"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'4.799"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'0.6291"
$ws.Range("E14").Value = "  -6.71%  "
$ws.Range("D15").Value = "'0.000009896"
$ws.Range("E15").Value = "  -3.66%  "
$ws.Range("D16").Value = "2.084.59"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'79.12"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'5.856"
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "29.269.15"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'224.53"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'11.71"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").Value = "'7.015"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'156.04"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").Value = "'0.1308"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'7.993"
$ws.Range("E27").Value = "  -5.91%  "
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "'0.06370"
$ws.Range("E30").Value = "  -12.62%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'3.849"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").Value = "'3.799"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("D34").Value = "'1.111"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "'1.732"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").Value = "'0.6474"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").Value = "'2.547"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "1.218.25"
$ws.Range("D39").Value = "'2.742"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").Value = "'0.01754"
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("D41").Value = "'6.549"
$ws.Range("E41").Value = "  -5.24%  "
$ws.Range("D42").Value = "'0.9004"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "1.997.06"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'100.19"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "'62.75"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "'0.00000000115"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("D48").Value = "'8.557"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").Value = "'1.577"
$ws.Range("E49").Value = "  -8.03%  "
$ws.Range("D50").Value = "'0.4557"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -2.63%  "
